# Update "want to go" counts (column F) on both the "展览" sheet and the
# mirrored "全部类型" sheet, which hold identical data in this workbook.
$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 238
    4  = 13151
    12 = 5
    18 = 5565
    22 = 17
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
